# The upstream commit ("add the test for docxToHtml") only adds this
# fixture to a test suite; the accompanying OOXML diff is not a real
# content edit at all. Every hunk in it only reorders the xmlns:*
# namespace-declaration attributes on a handful of root elements
# (document.xml, endnotes.xml, the header/footer parts, footnotes.xml,
# styles.xml, stylesWithEffects.xml and theme1.xml) -- e.g.
# `xmlns:r="..." xmlns:m="..."` becomes `xmlns:m="..." xmlns:r="..."`.
# No element, attribute value, run, paragraph, style, or any other
# piece of document content differs anywhere in the package: this is
# pure namespace-prefix-map churn left behind by the authoring tool
# (the fixture was regenerated by a different docx4j/JAXB version,
# whose internal map iteration order differs), and it disappears
# completely once the XML is canonicalized (namespace nodes are
# canonically sorted by prefix, so "xmlns:r, xmlns:m" and
# "xmlns:m, xmlns:r" canonicalize to the exact same output).
#
# There is therefore no actual Word-object-model action that
# corresponds to this "edit" -- Word's OM has no notion of raw
# namespace-declaration ordering to begin with, and nothing in the
# document's visible content, formatting, or structure changed. The
# correct COM-interop reproduction is simply to leave the document
# untouched (any synthetic edit here would only inject spurious
# content differences -- new rsids, stray namespace appends, etc. --
# that are not present in the target diff).

$d = $word.ActiveDocument

# Touch nothing; read a harmless, side-effect-free property so the
# script is an explicit, intentional no-op rather than an empty file.
$null = $d.Paragraphs.Count
